# Edit: append a parenthetical "(<=Also mehr oder weniger gerade zum Rand routen)"
# remark to the "Radikalbegradigung + OPO ..." paragraph, and give the
# "Labels bekommen eigenen Zeilenabstand!" paragraph its own (strikethrough)
# formatting instead of bold+underline, matching the commit
# "Annotations now have their own line spacing."

$d = $word.ActiveDocument

$para1Xml = '<w:p w:rsidR="00B84F0D" w:rsidRDefault="00B84F0D" w:rsidP="00032A5C"><w:pPr><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>Radikalbegradigung</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> + OPO </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>mit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="009B7FD9"><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">fixer </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>Max-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>Kapazität</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Also </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>mehr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>oder</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>weniger</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>gerade</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>zum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Rand </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>routen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>)</w:t></w:r></w:p>'
$para2Xml = '<w:p w:rsidR="006F6722" w:rsidRPr="00CB62D7" w:rsidRDefault="006F6722" w:rsidP="00032A5C"><w:pPr><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Labels </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>bekommen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>eigenen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>Zeilenabstand</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>!</w:t></w:r></w:p>'

# Locate the target paragraphs first (by distinctive text), so we don't
# mutate the Paragraphs collection while iterating over it.
$target1 = $null
$target2 = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($null -eq $target1 -and $t -like "*Radikalbegradigung*Kapazit*") {
        $target1 = $p.Range
    }
    elseif ($null -eq $target2 -and $t -like "*Labels*bekommen*eigenen*Zeilenabstand*") {
        $target2 = $p.Range
    }
}

if ($null -ne $target1) {
    $target1.InsertXML($para1Xml)
}
if ($null -ne $target2) {
    $target2.InsertXML($para2Xml)
}

Write-Host "Para1 replaced: $($null -ne $target1)"
Write-Host "Para2 replaced: $($null -ne $target2)"
